# Auto-refresh of cryptos list data (prices & 1h volume deltas).
# Mirrors the upstream GitHub Actions scrape-and-commit job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some refreshed prices are plain decimal-looking strings (e.g. "140.25").
# The sheet stores every Price/Volume cell as literal text (not a number),
# so mark those specific cells as Text-formatted first; this stops Excel's
# COM layer from auto-coercing the assigned string into a numeric value.
$forceTextCells = @("D5","D6","D8","D10","D12","D16","D18","D19","D22","D24","D26","D32","D34","D37","D38","D39","D40","D42","D43","D44","D45","D46","D47","D48","D50")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row-by-row cell updates (Coin / Link / Price / Volume(1h)).
# Row 2
$ws.Range("D2").Value = '59.250.54'
$ws.Range("E2").Value = '  +3.52%  '

# Row 3
$ws.Range("D3").Value = '2.593.26'
$ws.Range("E3").Value = '  +2.16%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '521.38'
$ws.Range("E5").Value = '  +1.56%  '

# Row 6
$ws.Range("D6").Value = '140.25'
$ws.Range("E6").Value = '  +0.97%  '

# Row 7
$ws.Range("E7").Value = '  -0.43%  '

# Row 8
$ws.Range("D8").Value = '0.566'
$ws.Range("E8").Value = '  +1.56%  '

# Row 9
$ws.Range("D9").Value = '2.616.81'
$ws.Range("E9").Value = '  +3.04%  '

# Row 10
$ws.Range("D10").Value = '6.52'
$ws.Range("E10").Value = '  +0.54%  '

# Row 11
$ws.Range("E11").Value = '  +1.94%  '

# Row 12
$ws.Range("D12").Value = '0.331'
$ws.Range("E12").Value = '  +2.09%  '

# Row 13
$ws.Range("E13").Value = '  +1.98%  '

# Row 14
$ws.Range("D14").Value = '3.053.96'
$ws.Range("E14").Value = '  +2.19%  '

# Row 15
$ws.Range("D15").Value = '59.199.62'
$ws.Range("E15").Value = '  +3.41%  '

# Row 16
$ws.Range("D16").Value = '20.42'
$ws.Range("E16").Value = '  +2.24%  '

# Row 17
$ws.Range("D17").Value = '2.607.10'
$ws.Range("E17").Value = '  +0.86%  '

# Row 18
$ws.Range("D18").Value = '0.0000132'
$ws.Range("E18").Value = '  +0.25%  '

# Row 19
$ws.Range("D19").Value = '337.91'
$ws.Range("E19").Value = '  +1.56%  '

# Row 20
$ws.Range("E20").Value = '  +1.18%  '

# Row 21
$ws.Range("E21").Value = '  +1.20%  '

# Row 22
$ws.Range("D22").Value = '6.49'
$ws.Range("E22").Value = '  +5.98%  '

# Row 23
$ws.Range("E23").Value = '  -0.23%  '

# Row 24
$ws.Range("D24").Value = '66.35'
$ws.Range("E24").Value = '  +2.79%  '

# Row 25
$ws.Range("E25").Value = '  +1.02%  '

# Row 26
$ws.Range("D26").Value = '0.403'
$ws.Range("E26").Value = '  +0.89%  '

# Row 27
$ws.Range("E27").Value = '  -0.65%  '

# Row 28
$ws.Range("E28").Value = '  +1.65%  '

# Row 29
$ws.Range("E29").Value = '  -0.13%  '

# Row 30
$ws.Range("D30").Value = '0.0₃0725'
$ws.Range("E30").Value = '  -3.39%  '

# Row 31
$ws.Range("E31").Value = '  -4.68%  '

# Row 32
$ws.Range("D32").Value = '18.81'
$ws.Range("E32").Value = '  +1.94%  '

# Row 33
$ws.Range("E33").Value = '  +0.91%  '

# Row 34
$ws.Range("D34").Value = '149.16'
$ws.Range("E34").Value = '  +0.22%  '

# Row 35
$ws.Range("E35").Value = '  +0.98%  '

# Row 36
$ws.Range("E36").Value = '  +0.12%  '

# Row 37
$ws.Range("D37").Value = '36.35'
$ws.Range("E37").Value = '  +1.70%  '

# Row 38
$ws.Range("D38").Value = '1.45'
$ws.Range("E38").Value = '  +3.74%  '

# Row 39
$ws.Range("D39").Value = '0.832'
$ws.Range("E39").Value = '  +1.29%  '

# Row 40
$ws.Range("D40").Value = '0.823'
$ws.Range("E40").Value = '  -1.68%  '

# Row 41
$ws.Range("E41").Value = '  +1.91%  '

# Row 42
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = '0.996'
$ws.Range("E42").Value = '  -0.42%  '

# Row 43
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Value = '275.36'
$ws.Range("E43").Value = '  +6.98%  '

# Row 44
$ws.Range("D44").Value = '10.72'
$ws.Range("E44").Value = '  +0.96%  '

# Row 45
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '0.591'
$ws.Range("E45").Value = '  +3.04%  '

# Row 46
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = '0.0953'
$ws.Range("E46").Value = '  +0.12%  '

# Row 47
$ws.Range("D47").Value = '0.0520'
$ws.Range("E47").Value = '  +0.09%  '

# Row 48
$ws.Range("D48").Value = '18.56'
$ws.Range("E48").Value = '  +0.88%  '

# Row 49
$ws.Range("D49").Value = '1.981.02'
$ws.Range("E49").Value = '  +0.81%  '

# Row 50
$ws.Range("D50").Value = '4.61'
$ws.Range("E50").Value = '  +2.73%  '

# Row 51
$ws.Range("E51").Value = '  -0.21%  '
